$wb = $excel.ActiveWorkbook

$ws2010 = $wb.Worksheets.Item("2010")
$ws2010.Range("K7").ClearContents()
$ws2010.Range("L7").ClearContents()
